{"js": "// Integra\u00e7\u00e3o dos inputs do formul\u00e1rio com o dicion\u00e1rio de convers\u00e3o\n//\n// Replaces the declarant's name (\"Marcia\" -> \"dfgbdsfgbdfgb\", both\n// occurrences: in the body paragraph and in the signature-block table),\n// replaces the personal/address details run, and updates the date.\n\nconst body = context.document.body;\n\n// 1) Name \"Marcia\" -> \"dfgbdsfgbdfgb\" (occurs twice: body paragraph + table cell)\nconst nameResults = body.search(\"Marcia\", { matchCase: true, matchWholeWord: true });\nnameResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < nameResults.items.length; i++) {\n  nameResults.items[i].insertText(\"dfgbdsfgbdfgb\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Qualification / address run\nconst oldDetails =\n  \", brasileira, solteira, professora, registrada no CPF sob o n\u00ba 058.252.897-63, residente e domiciliada na Rua Dona Romana, n\u00ba 406, apto. 102, Engenho Novo, Rio de Janeiro - RJ, CEP 20710-200, declara, para os devidos fins, que n\u00e3o possui condi\u00e7\u00f5es financeiras para arcar com as custas processuais e honor\u00e1rios advocat\u00edcios da presente demanda sem preju\u00edzo do pr\u00f3prio sustento e o de sua fam\u00edlia.\";\nconst newDetails =\n  \", Brasileiro, Casado, fdxgeg, registrada no CPF sob o n\u00ba 000000000, residente e domiciliada na Estrada fgdfghdsgh, 0555, xcvgfdsg, dfgfdgdfg, dfgdfrgd - CE, CEP: 000000, declara, para os devidos fins, que n\u00e3o possui condi\u00e7\u00f5es financeiras para arcar com as custas processuais e honor\u00e1rios advocat\u00edcios da presente demanda sem preju\u00edzo do pr\u00f3prio sustento e o de sua fam\u00edlia.\";\n\nconst detailsResults = body.search(oldDetails, { matchCase: true });\ndetailsResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < detailsResults.items.length; i++) {\n  detailsResults.items[i].insertText(newDetails, \"Replace\");\n}\nawait context.sync();\n\n// 3) Date \"23 de mar\u00e7o de 2024\" -> \"25 de mar\u00e7o de 2024\"\nconst dateResults = body.search(\"23 de mar\u00e7o de 2024\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"25 de mar\u00e7o de 2024\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Integra\u00e7\u00e3o dos inputs do formul\u00e1rio com o dicion\u00e1rio de convers\u00e3o\n#\n# Replaces the declarant's name (\"Marcia\" -> \"dfgbdsfgbdfgb\", both\n# occurrences: in the body paragraph and in the signature-block table),\n# replaces the personal/address details run, and updates the date.\n#\n# Each replacement is scoped through a temporary Bookmark around the\n# Find hit before assigning Range.Text; this performs a surgical,\n# single-run rewrite instead of letting Word's Find/Replace collapse\n# neighbouring same-formatted runs in the paragraph into one run.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllExact($doc, $oldText, $newText, $bookmarkBase) {\n    $count = 0\n    while ($true) {\n        $find = $doc.Content.Find\n        $find.ClearFormatting()\n        $find.Text = $oldText\n        $find.Forward = $true\n        $find.Wrap = 0\n        $found = $find.Execute()\n        if (-not $found) { break }\n\n        $r = $find.Parent\n        $bmName = \"$bookmarkBase$count\"\n        $doc.Bookmarks.Add($bmName, $r) | Out-Null\n        $bmRange = $doc.Bookmarks.Item($bmName).Range\n        $bmRange.Text = $newText\n        $doc.Bookmarks.Item($bmName).Delete()\n\n        $count = $count + 1\n        if ($count -gt 50) { break }\n    }\n    return $count\n}\n\n# 1) Name \"Marcia\" -> \"dfgbdsfgbdfgb\" (replace all occurrences: body + table)\nReplace-AllExact $d \"Marcia\" \"dfgbdsfgbdfgb\" \"tmpName\" | Out-Null\n\n# 2) Qualification / address run\n$oldDetails = \", brasileira, solteira, professora, registrada no CPF sob o n\u00ba 058.252.897-63, residente e domiciliada na Rua Dona Romana, n\u00ba 406, apto. 102, Engenho Novo, Rio de Janeiro - RJ, CEP 20710-200, declara, para os devidos fins, que n\u00e3o possui condi\u00e7\u00f5es financeiras para arcar com as custas processuais e honor\u00e1rios advocat\u00edcios da presente demanda sem preju\u00edzo do pr\u00f3prio sustento e o de sua fam\u00edlia.\"\n$newDetails = \", Brasileiro, Casado, fdxgeg, registrada no CPF sob o n\u00ba 000000000, residente e domiciliada na Estrada fgdfghdsgh, 0555, xcvgfdsg, dfgfdgdfg, dfgdfrgd - CE, CEP: 000000, declara, para os devidos fins, que n\u00e3o possui condi\u00e7\u00f5es financeiras para arcar com as custas processuais e honor\u00e1rios advocat\u00edcios da presente demanda sem preju\u00edzo do pr\u00f3prio sustento e o de sua fam\u00edlia.\"\nReplace-AllExact $d $oldDetails $newDetails \"tmpDetails\" | Out-Null\n\n# 3) Date \"23 de mar\u00e7o de 2024\" -> \"25 de mar\u00e7o de 2024\"\nReplace-AllExact $d \"23 de mar\u00e7o de 2024\" \"25 de mar\u00e7o de 2024\" \"tmpDate\" | Out-Null\n\nWrite-Output \"done\"\n"}
